# Apply the "gh-pages output generated at 456a3b4" update:
#  - bump several "want-to-go" counters across sheets
#  - on "展览" and "全部类型" insert a brand-new event
#    ("张家港·授渔4.0圆梦展-阿判的超绝二次元大趴") right before the existing
#    "苏州·Good jump ACG元旦跨年盛典国潮文化节" row, pushing it (and the row
#    after it) down by one, and bump their counters too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value = 175
$ws1.Range("F6").Value = 22
$ws1.Range("F8").Value = 393
$ws1.Range("F9").Value = 1011
$ws1.Range("F12").Value = 536
$ws1.Range("F14").Value = 12637
$ws1.Range("G14").Value = 49.9

# Row 15 currently holds "Good jump ACG...". Its A/B columns (14 /
# 2025-01-01) stay put; insert a fresh blank row *below* it (row 16) to
# hold the old "Good jump ACG..." row contents (with an updated
# want-to-go count), which also pushes the old row 16 ("星部落...") down
# to row 17.
$ws1.Rows.Item(16).Insert()
$ws1.Range("A15").Copy()
$ws1.Range("A16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# New row 16: what used to be in row 15 ("Good jump ACG..."), with an
# updated want-to-go count.
$ws1.Range("A16").Value = 15
$ws1.Range("B16").Value = "2025-01-01"
$ws1.Range("C16").Value = "苏州·Good jump ACG元旦跨年盛典国潮文化节"
$ws1.Range("D16").Value = "金山南路影视城 木渎影视城会展中心"
$ws1.Range("E16").Value = "2025.01.01 10:00-01.01 17:00"
$ws1.Range("F16").Value = 5211
$ws1.Range("G16").Value = 60
$ws1.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=93234"
$ws1.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202410/Aq3TKqhc1728483225862.jpeg"

# Row 15 is overwritten in place with the brand-new event (A15/B15 keep
# their existing values: 14 / 2025-01-01).
$ws1.Range("C15").Value = "张家港·授渔4.0圆梦展-阿判的超绝二次元大趴"
$ws1.Range("D15").Value = "泗杨路张家港碧桂园天玺东南侧约60米 五月风华宴会中心"
$ws1.Range("E15").Value = "2025.01.01 09:30-01.01 16:30"
$ws1.Range("F15").Value = 0
$ws1.Range("G15").Value = 40
$ws1.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=93391"
$ws1.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202410/59vdXP5U1728832042854.png"

# Row 17 now holds what used to be row 16 ("星部落...") with an updated
# want-to-go count.
$ws1.Range("F17").Value = 5519

# ---------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 42

# ---------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F6").Value = 175
$ws4.Range("F7").Value = 22
$ws4.Range("F9").Value = 393
$ws4.Range("F10").Value = 1011
$ws4.Range("F13").Value = 536
$ws4.Range("F15").Value = 12637
$ws4.Range("G15").Value = 49.9
$ws4.Range("F16").Value = 42

# Row 18 currently holds "Good jump ACG...". Its A/B columns (17 /
# 2025-01-01) stay put; insert a fresh blank row *below* it (row 19) to
# hold the old "Good jump ACG..." row contents (with an updated
# want-to-go count), which also pushes the old row 19 ("星部落...") down
# to row 20.
$ws4.Rows.Item(19).Insert()
$ws4.Range("A18").Copy()
$ws4.Range("A19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# New row 19: what used to be in row 18 ("Good jump ACG..."), with an
# updated want-to-go count.
$ws4.Range("A19").Value = 18
$ws4.Range("B19").Value = "2025-01-01"
$ws4.Range("C19").Value = "苏州·Good jump ACG元旦跨年盛典国潮文化节"
$ws4.Range("D19").Value = "金山南路影视城 木渎影视城会展中心"
$ws4.Range("E19").Value = "2025.01.01 10:00-01.01 17:00"
$ws4.Range("F19").Value = 5211
$ws4.Range("G19").Value = 60
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=93234"
$ws4.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202410/Aq3TKqhc1728483225862.jpeg"

# Row 18 is overwritten in place with the brand-new event (A18/B18 keep
# their existing values: 17 / 2025-01-01).
$ws4.Range("C18").Value = "张家港·授渔4.0圆梦展-阿判的超绝二次元大趴"
$ws4.Range("D18").Value = "泗杨路张家港碧桂园天玺东南侧约60米 五月风华宴会中心"
$ws4.Range("E18").Value = "2025.01.01 09:30-01.01 16:30"
$ws4.Range("F18").Value = 0
$ws4.Range("G18").Value = 40
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=93391"
$ws4.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202410/59vdXP5U1728832042854.png"

# Row 20 now holds what used to be row 19 ("星部落...") with an updated
# want-to-go count.
$ws4.Range("F20").Value = 5519
